# Weekly update: insert the newest week's two records (Primera, 36-unit box and
# Primera, 60-unit box, both dated 45021 = 2023-04-05) at the top of the data
# block (right after the existing header/first-data row), pushing every other
# record down by two rows. This mirrors how the source weekly consolidation
# sheet keeps newest-first ordering for this market/product subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 666:667 (existing rows 666.. shift down to 668..)
$ws.Rows("666:667").Insert()

# --- New row 666 ---
$ws.Cells.Item(666, 1).Value  = 3
$ws.Cells.Item(666, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(666, 3).Value  = "Coquimbo"
$ws.Cells.Item(666, 4).Value  = 45021
$ws.Cells.Item(666, 5).Value  = 5
$ws.Cells.Item(666, 6).Value  = 100112032
$ws.Cells.Item(666, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(666, 8).Value  = "Sin especificar"
$ws.Cells.Item(666, 9).Value  = "Primera"
$ws.Cells.Item(666, 10).Value = 160
$ws.Cells.Item(666, 11).Value = 4000
$ws.Cells.Item(666, 12).Value = 4500
$ws.Cells.Item(666, 13).Value = 4250
$ws.Cells.Item(666, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(666, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(666, 16).Value = 118
$ws.Cells.Item(666, 17).Value = 36
$ws.Cells.Item(666, 18).Value = "Hortaliza"

# --- New row 667 ---
$ws.Cells.Item(667, 1).Value  = 3
$ws.Cells.Item(667, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(667, 3).Value  = "Coquimbo"
$ws.Cells.Item(667, 4).Value  = 45021
$ws.Cells.Item(667, 5).Value  = 5
$ws.Cells.Item(667, 6).Value  = 100112032
$ws.Cells.Item(667, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(667, 8).Value  = "Sin especificar"
$ws.Cells.Item(667, 9).Value  = "Primera"
$ws.Cells.Item(667, 10).Value = 220
$ws.Cells.Item(667, 11).Value = 7000
$ws.Cells.Item(667, 12).Value = 7500
$ws.Cells.Item(667, 13).Value = 7266
$ws.Cells.Item(667, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(667, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(667, 16).Value = 121
$ws.Cells.Item(667, 17).Value = 60
$ws.Cells.Item(667, 18).Value = "Hortaliza"
